# Update the "Team members" slide (sldId 259 -> 3rd slide in the deck) so
# that the 3rd bullet reads "Shriyansh Akash Jain [2204030100404]" instead
# of "Nishith Narendra Mehta [2204030102026]".
#
# The paragraph originally has 2 runs:
#   1) "Nishith Narendra Mehta "   (lang="en-IN")
#   2) "[2204030102026]"          (lang="en-US")
#
# It needs to become 3 runs:
#   1) "Shriyansh"                (lang="en-IN", newly inserted)
#   2) " Akash Jain "             (lang="en-IN", same run/formatting as old run 1)
#   3) "[2204030100404]"          (lang="en-US", same run/formatting as old run 2)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 2") {
        $target = $sh
        break
    }
}

$tr = $target.TextFrame.TextRange

# This textbox has <a:spAutoFit/>; editing its text makes the host
# recompute the box height. Capture the original height so it can be
# restored afterwards (the source edit did not change the box size).
$origHeight = $target.Height

# Locate the paragraph that holds the "Nishith Narendra Mehta" entry.
$paraCount = $tr.Paragraphs().Count
$paraIndex = 0
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.IndexOf("Nishith Narendra Mehta") -ge 0) {
        $paraIndex = $i
        break
    }
}

$para = $tr.Paragraphs($paraIndex, 1)

# Insert the new "Shriyansh" run immediately before the paragraph's
# existing text; it inherits the formatting of the first existing run.
$para.InsertBefore("Shriyansh") | Out-Null

# Replace the old name text (now shifted right by len("Shriyansh")) with
# " Akash Jain ", keeping the same run/formatting.
$full = $tr.Text
$nameStart = $full.IndexOf("Nishith Narendra Mehta ")
$nameRange = $tr.Characters($nameStart + 1, 23)
$nameRange.Text = " Akash Jain "

# Replace the old roll number text with the new one, keeping the same
# run/formatting.
$full = $tr.Text
$idStart = $full.IndexOf("[2204030102026]")
$idRange = $tr.Characters($idStart + 1, 15)
$idRange.Text = "[2204030100404]"

# Restore the original autofit height (nudged by half an EMU-in-points so
# the point->EMU rounding lands back on the original integer EMU value).
$target.Height = $origHeight + 0.00004

